# "cambios en el servidor" - a new guest/RSVP entry was appended to the
# guest list sheet (row 2): a test submission with name "PRUEBA CEL",
# sender "Novio" and a congratulations message "PRIEBA CEL \n  ...".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row of data right below the header row.
$ws.Cells.Item(2, 1).Value = "PRUEBA CEL"
$ws.Cells.Item(2, 2).Value = "Novio"
$ws.Cells.Item(2, 3).Value = "PRIEBA CEL `n                    "

# The congratulations message contains an embedded newline, which would
# otherwise make the engine auto-compute a taller row height for row 2.
# Reset it back to the sheet's standard height and autofit so the row
# keeps its default (non-custom) height, matching a simple data append.
$ws.Rows.Item(2).RowHeight = 15
$ws.Rows.Item(2).EntireRow.AutoFit()
